$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "transmit, receive, idle energy model" values for column A (LEACH)
$ws.Range("A2").Value = 8
$ws.Range("A3").Value = 2
$ws.Range("A4:A5").Value = 8
$ws.Range("A7:A15").Value = 8
$ws.Range("A16:A17").Value = 2
$ws.Range("A19").Value = 4
$ws.Range("A22:A202").Value = 1
